$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Kerabat Mempelai Pria"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Keluarga Mempelai Wanita"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Keluarga Mempelai Wanita"),
    @("DDMMYYFN20", "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5, "Keluarga Mempelai Wanita")
)

# Column D holds a phone number with a leading zero; format as Text first
# so Excel doesn't silently coerce it to a number and drop the leading zero.
$ws.Range("D12:D22").NumberFormat = "@"

$startRow = 12
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
